# Reorder the "Recorded By" (column G) comma-separated list of recorders
# so that whenever the first entry in the list is "System" (case-insensitive),
# the whole list order is reversed. Rows where "System" is not the first
# entry (or the cell has only a single value) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $value = $cell.Value2

    if ($value -eq $null) { continue }

    $text = [string]$value
    if ($text -eq "") { continue }

    $parts = $text -split ",\s*"
    if ($parts.Count -le 1) { continue }

    if ($parts[0].ToLower() -eq "system") {
        $reversedParts = @()
        for ($i = $parts.Count - 1; $i -ge 0; $i--) {
            $reversedParts += $parts[$i]
        }
        $newText = [string]::Join(", ", $reversedParts)
        $cell.Value = $newText
    }
}
